# Generate Report for Handoff
# Updates the localization-status workbook: the source file was renamed
# (new GUID) and re-handed-off, so the "target"/"handback" columns for
# the not-yet-handed-back languages are cleared again.

$wb = $excel.ActiveWorkbook

$oldGuid = "bc7ad01a-d8c4-4593-bace-17fb2811f112"
$newGuid = "74acc957-0367-4cb2-b2fa-dbab47ec869f"
$oldHash = "531a9394bd1e7a4793c0429ba8ba9aa0cc169170"
$newHash = "440fb37ad80315ebdc302b739831cb373bfbc829"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# File Name (A2) carries the new GUID too.
$wsOverview.Range("A2").Value = $newGuid + ".md"

# Path And Name (B2) is a hyperlink; both the cell text and the display
# text change (same string) since the source file's GUID is new.
$wsOverview.Range("B2").Value = "e2e\" + $newGuid + ".md"
foreach ($h in $wsOverview.Hyperlinks) {
    if ($h.Range.Row -eq 2 -and $h.Range.Column -eq 2) {
        $h.TextToDisplay = "e2e\" + $newGuid + ".md"
    }
}

# Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = "2016-08-13 19:15:50"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Source File Name (A2) is a hyperlink; update both cell text and display text.
$wsZh.Range("A2").Value = $newGuid + ".md"
foreach ($h in $wsZh.Hyperlinks) {
    if ($h.Range.Row -eq 2 -and $h.Range.Column -eq 1) {
        $h.TextToDisplay = $newGuid + ".md"
    }
}

# Latest Handoff File / Datetime
$wsZh.Range("G2").Value = $newGuid + "." + $newHash + ".zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-13 19:15:42"

# Latest Target File (I2) had a hyperlink; it's cleared out entirely
# (no handback yet against the new handoff), so drop the hyperlink too.
foreach ($h in $wsZh.Hyperlinks) {
    if ($h.Range.Row -eq 2 -and $h.Range.Column -eq 9) {
        $h.Delete()
    }
}
$wsZh.Range("I2").Value = ""

# Latest Handback File is likewise cleared.
$wsZh.Range("J2").Value = ""

# Latest Handback DateTime resets to the "never" sentinel.
$wsZh.Range("K2").Value = "0001-01-01 00:00:00"

# Latest Target File / Latest Handback File columns are narrower now that
# they hold short/empty values instead of full file names.
$wsZh.Columns.Item(9).ColumnWidth = 17.8333333333333
$wsZh.Columns.Item(10).ColumnWidth = 20.8333333333333

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newGuid + ".md"
foreach ($h in $wsDe.Hyperlinks) {
    if ($h.Range.Row -eq 2 -and $h.Range.Column -eq 1) {
        $h.TextToDisplay = $newGuid + ".md"
    }
}

# Latest Handoff File / Datetime (the Handoff Datetime shares the same
# text as Overview's Latest HO Xliff Generate Date).
$wsDe.Range("G2").Value = $newGuid + "." + $newHash + ".de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-13 19:15:50"

foreach ($h in $wsDe.Hyperlinks) {
    if ($h.Range.Row -eq 2 -and $h.Range.Column -eq 9) {
        $h.Delete()
    }
}
$wsDe.Range("I2").Value = ""
$wsDe.Range("J2").Value = ""
$wsDe.Range("K2").Value = "0001-01-01 00:00:00"

$wsDe.Columns.Item(9).ColumnWidth = 17.8333333333333
$wsDe.Columns.Item(10).ColumnWidth = 20.8333333333333
